# "updating excel loot values"
# Fill in the "Amount" column (B) on the "Loot" sheet with the loot values
# for each lootable object/character listed in column A, and leave the
# sheet's selection on that newly-filled range (matching the on-screen
# state after typing the values in).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Loot")
$ws.Activate()

# Row 2 -> Row 11 correspond to: Sapphire, Emerald, Coin, Coin Pile,
# Wad of Cash, Gold Bars, Wallet, Cops, Classy Woman, Timmy
$lootValues = @(25, 10, 15, 5, 250, 250, 15, 100, 200, 50)

for ($i = 0; $i -lt $lootValues.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $lootValues[$i]
    $cell.Font.Size = 16
}

$ws.Range("B2:B11").Select()
